# Regenerate save_data to use K instead of Strike#.
# This updates the "K" column (column G) values for rows 2-53 in place,
# writing the newly-calculated s_vals (K counts) that replace the old
# Strike# derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,1,3,0,0,0,2,0,0,0,2,1,0,1,0,0,0,1,2,0,1,1,0,3,1,2,1,4,1,0,0,0,1,0,1,2,0,0,0,0,1,1,1,1,0,0,1,1,1,0,0,0)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
